$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores plain decimal-looking strings (e.g. "1.003").
# Excel would otherwise auto-convert these to numbers on assignment, so
# each such cell is switched to Text format first, then restored to the
# default style afterwards (matching the original inline-string cells).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.556.11"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "1.746.82"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "321.46"
$ws.Range("E5").Value = "  +0.09%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "0.4614"
$ws.Range("E7").Value = "  +9.10%  "
$ws.Range("D8").Value = "0.3540"
$ws.Range("E8").Value = "  -2.28%  "
$ws.Range("D9").Value = "0.07444"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "42.00"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "1.089"
$ws.Range("E11").Value = "  +0.38%  "
$ws.Range("D12").Value = "1.003"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "20.61"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "5.961"
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "7.098"
$ws.Range("D16").Value = "1.742.85"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "91.66"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").Value = "0.00001058"
$ws.Range("E18").Value = "  +0.56%  "
$ws.Range("D19").Value = "0.06403"
$ws.Range("E19").Value = "  +0.84%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "16.69"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").Value = "5.753"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "27.618.51"
$ws.Range("E23").Value = "  -0.63%  "
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  -0.37%  "
$ws.Range("E25").Value = "  +0.35%  "
$ws.Range("D26").Value = "162.87"
$ws.Range("E26").Value = "  +3.51%  "
$ws.Range("D27").Value = "20.05"
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("D29").Value = "125.51"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("D30").Value = "2.049"
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("D31").Value = "1.048"
$ws.Range("E31").Value = "  -6.26%  "
$ws.Range("D32").Value = "0.09239"
$ws.Range("E32").Value = "  +4.80%  "
$ws.Range("D33").Value = "3.669"
$ws.Range("E33").Value = "  +0.07%  "
$ws.Range("D34").Value = "5.495"
$ws.Range("E34").Value = "  -0.96%  "
$ws.Range("D35").Value = "0.02281"
$ws.Range("E35").Value = "  -0.17%  "
$ws.Range("D36").Value = "11.71"
$ws.Range("E36").Value = "  -4.47%  "
$ws.Range("D37").Value = "0.06022"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "0.2079"
$ws.Range("E38").Value = "  -0.49%  "
$ws.Range("D39").Value = "4.935"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").Value = "0.6267"
$ws.Range("E40").Value = "  -0.30%  "
$ws.Range("D41").Value = "1.190"
$ws.Range("E41").Value = "  +1.16%  "
$ws.Range("D42").Value = "1.378"
$ws.Range("E42").Value = "  -1.56%  "
$ws.Range("D43").Value = "7.710"
$ws.Range("E43").Value = "  -1.62%  "
$ws.Range("D44").Value = "13.18"
$ws.Range("E44").Value = "  -1.19%  "
$ws.Range("D45").Value = "3.697"
$ws.Range("E45").Value = "  +0.58%  "
$ws.Range("D46").Value = "0.5844"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "122.24"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "1.931"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.06859"
$ws.Range("E49").Value = "  +0.80%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "1.130"
$ws.Range("E50").Value = "  -3.91%  "
$ws.Range("D51").Value = "71.46"
$ws.Range("E51").Value = "  -2.66%  "

$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
